$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to be treated as text so that numeric-looking
    # strings (e.g. '0.999', '0.0000186') are preserved verbatim,
    # then restore the default 'Normal' style so no stray number
    # format is left behind on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '64.042.86'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '2.650.47'
$ws.Range("E3").Value = '  +0.73%  '
Set-TextValue D4 '0.999'
$ws.Range("E4").Value = '  -0.12%  '
Set-TextValue D5 '582.45'
$ws.Range("E5").Value = '  +0.19%  '
Set-TextValue D6 '156.72'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -3.52%  '
$ws.Range("D9").Value = '2.645.08'
$ws.Range("E9").Value = '  +0.59%  '
Set-TextValue D10 '0.119'
$ws.Range("E10").Value = '  -2.96%  '
Set-TextValue D11 '5.82'
$ws.Range("E11").Value = '  +0.20%  '
Set-TextValue D12 '0.384'
$ws.Range("E12").Value = '  -1.41%  '
Set-TextValue D13 '0.157'
$ws.Range("E13").Value = '  +1.09%  '
Set-TextValue D14 '28.68'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '3.125.10'
$ws.Range("E15").Value = '  +0.53%  '
Set-TextValue D16 '0.0000186'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = '63.924.83'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '2.641.92'
$ws.Range("E18").Value = '  +0.19%  '
Set-TextValue D19 '12.24'
$ws.Range("E19").Value = '  +0.04%  '
Set-TextValue D20 '7.77'
$ws.Range("E20").Value = '  +4.66%  '
$ws.Range("E21").Value = '  -2.61%  '
Set-TextValue D22 '346.57'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +0.36%  '
Set-TextValue D24 '68.06'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  +5.15%  '
Set-TextValue D26 '0.0000113'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue D27 '9.33'
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue D28 '595.81'
$ws.Range("E28").Value = '  +1.66%  '
Set-TextValue D29 '1.62'
$ws.Range("E29").Value = '  +2.45%  '
Set-TextValue D30 '8.24'
$ws.Range("E30").Value = '  +3.56%  '
$ws.Range("E31").Value = '  +0.38%  '
Set-TextValue D32 '0.999'
$ws.Range("E32").Value = '  -0.24%  '
Set-TextValue D33 '2.08'
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("E34").Value = '  +1.24%  '
Set-TextValue D35 '6.67'
$ws.Range("E35").Value = '  +0.28%  '
Set-TextValue D36 '5.52'
$ws.Range("E36").Value = '  +3.63%  '
Set-TextValue D37 '0.405'
$ws.Range("E37").Value = '  -1.77%  '
Set-TextValue D38 '19.79'
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("E39").Value = '  -0.07%  '
Set-TextValue D40 '1.93'
$ws.Range("E40").Value = '  +0.02%  '
Set-TextValue D41 '150.83'
$ws.Range("E41").Value = '  -2.43%  '
Set-TextValue D42 '2.56'
$ws.Range("E42").Value = '  +4.78%  '
Set-TextValue D45 '163.92'
$ws.Range("E45").Value = '  +3.90%  '
Set-TextValue D46 '24.28'
$ws.Range("E46").Value = '  +4.69%  '
Set-TextValue D47 '3.92'
$ws.Range("E47").Value = '  -2.00%  '
Set-TextValue D48 '0.0592'
$ws.Range("E48").Value = '  -1.43%  '
Set-TextValue D49 '0.636'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("E51").Value = '  -1.73%  '
